$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")
$ws.Range("A7").Value = "A7"
